$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add daily power records for rows 53-56 (following the existing pattern of
# date / start time / end time / duration formulas used throughout the
# table).
# ---------------------------------------------------------------------------

# Row 53 (2018-10-04): no start/end time recorded, only the date.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$ws.Range("D52:F52").Copy()
$ws.Range("D53:F53").PasteSpecial(-4122)

$ws.Range("A53").Value = 43377
$ws.Range("D53").Formula = "=(C53-B53)* 1440"
$ws.Range("E53").Formula = "=IF(C53>B53, (C53-B53)*1440, (B53-C53)*1440)"
$ws.Range("F53").Formula = "=ABS((C53-B53)*1440)"

# Row 54 (2018-10-05): full record with start and end time.
$ws.Range("A52:F52").Copy()
$ws.Range("A54:F54").PasteSpecial(-4122)

$ws.Range("A54").Value = 43378
$ws.Range("B54").Value = 0.77430555555555547
$ws.Range("C54").Value = 0.99930555555555556
$ws.Range("D54").Formula = "=(C54-B54)* 1440"
$ws.Range("E54").Formula = "=IF(C54>B54, (C54-B54)*1440, (B54-C54)*1440)"
$ws.Range("F54").Formula = "=ABS((C54-B54)*1440)"

# Row 55 (2018-10-06): full record with start and end time.
$ws.Range("A52:F52").Copy()
$ws.Range("A55:F55").PasteSpecial(-4122)

$ws.Range("A55").Value = 43379
$ws.Range("B55").Value = 0
$ws.Range("C55").Value = 0.31527777777777777
$ws.Range("D55").Formula = "=(C55-B55)* 1440"
$ws.Range("E55").Formula = "=IF(C55>B55, (C55-B55)*1440, (B55-C55)*1440)"
$ws.Range("F55").Formula = "=ABS((C55-B55)*1440)"

# Row 56 (2018-10-07): no start/end time recorded, only the date.
$ws.Range("A52").Copy()
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("D52:F52").Copy()
$ws.Range("D56:F56").PasteSpecial(-4122)

$ws.Range("A56").Value = 43380
$ws.Range("D56").Formula = "=(C56-B56)* 1440"
$ws.Range("E56").Formula = "=IF(C56>B56, (C56-B56)*1440, (B56-C56)*1440)"
$ws.Range("F56").Formula = "=ABS((C56-B56)*1440)"

# Extend the table range to cover the newly added rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F56"))

# Update the active selection to reflect the new last cell, matching the
# scrolled view after the records were appended.
[void]$ws.Range("A45").Select()
[void]$ws.Range("B56").Select()
